$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'57.435.72"
$ws.Range("E2").Value = "  -1.68%  "

# Row 3
$ws.Range("D3").Value = "'2.426.70"
$ws.Range("E3").Value = "  -2.27%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").Value = "'509.31"
$ws.Range("E5").Value = "  -3.00%  "

# Row 6
$ws.Range("D6").Value = "'127.95"
$ws.Range("E6").Value = "  -4.26%  "

# Row 7
$ws.Range("E7").Value = "  -0.10%  "

# Row 8
$ws.Range("E8").Value = "  -2.44%  "

# Row 9
$ws.Range("D9").Value = "'2.437.25"
$ws.Range("E9").Value = "  -1.89%  "

# Row 10
$ws.Range("E10").Value = "  -0.36%  "

# Row 11
$ws.Range("D11").Value = "'0.0944"
$ws.Range("E11").Value = "  -5.48%  "

# Row 12
$ws.Range("D12").Value = "'5.13"
$ws.Range("E12").Value = "  -5.49%  "

# Row 13
$ws.Range("D13").Value = "'0.328"
$ws.Range("E13").Value = "  -4.18%  "

# Row 14
$ws.Range("D14").Value = "'2.859.75"
$ws.Range("E14").Value = "  -2.25%  "

# Row 15
$ws.Range("D15").Value = "'57.352.38"
$ws.Range("E15").Value = "  -1.77%  "

# Row 16
$ws.Range("D16").Value = "'21.62"
$ws.Range("E16").Value = "  -3.58%  "

# Row 17
$ws.Range("E17").Value = "  -3.85%  "

# Row 18
$ws.Range("D18").Value = "'2.433.21"
$ws.Range("E18").Value = "  -2.26%  "

# Row 19
$ws.Range("D19").Value = "'10.37"
$ws.Range("E19").Value = "  -5.08%  "

# Row 20
$ws.Range("D20").Value = "'314.11"
$ws.Range("E20").Value = "  -2.30%  "

# Row 21
$ws.Range("D21").Value = "'4.08"
$ws.Range("E21").Value = "  -3.07%  "

# Row 22
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  +0.15%  "

# Row 23
$ws.Range("D23").Value = "'5.61"
$ws.Range("E23").Value = "  -3.47%  "

# Row 24
$ws.Range("D24").Value = "'63.14"
$ws.Range("E24").Value = "  -1.92%  "

# Row 25
$ws.Range("D25").Value = "'0.402"
$ws.Range("E25").Value = "  -2.60%  "

# Row 26
$ws.Range("D26").Value = "'0.997"
$ws.Range("E26").Value = "  -0.29%  "

# Row 27
$ws.Range("E27").Value = "  -1.95%  "

# Row 28
$ws.Range("D28").Value = "'7.19"
$ws.Range("E28").Value = "  -3.78%  "

# Row 29
$ws.Range("D29").Value = "'169.69"
$ws.Range("E29").Value = "  +2.02%  "

# Row 30
$ws.Range("E30").Value = "  -3.68%  "

# Row 31
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.66"
$ws.Range("E31").Value = "  -3.02%  "

# Row 32
$ws.Range("B32").Value = "PEPE"
$ws.Range("C32").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D32").Value = "'0.0₃0713"
$ws.Range("E32").Value = "  -5.17%  "

# Row 33
$ws.Range("D33").Value = "'1.15"
$ws.Range("E33").Value = "  +1.37%  "

# Row 34
$ws.Range("E34").Value = "  -0.05%  "

# Row 35
$ws.Range("D35").Value = "'0.998"
$ws.Range("E35").Value = "  -0.11%  "

# Row 36
$ws.Range("D36").Value = "'17.60"
$ws.Range("E36").Value = "  -3.52%  "

# Row 37
$ws.Range("E37").Value = "  -5.63%  "

# Row 38
$ws.Range("D38").Value = "'3.89"
$ws.Range("E38").Value = "  -2.26%  "

# Row 39
$ws.Range("D39").Value = "'36.18"
$ws.Range("E39").Value = "  -1.19%  "

# Row 40
$ws.Range("E40").Value = "  -3.69%  "

# Row 41
$ws.Range("D41").Value = "'0.765"
$ws.Range("E41").Value = "  -4.21%  "

# Row 42
$ws.Range("D42").Value = "'270.58"
$ws.Range("E42").Value = "  -2.78%  "

# Row 43
$ws.Range("D43").Value = "'3.35"
$ws.Range("E43").Value = "  -5.00%  "

# Row 44
$ws.Range("D44").Value = "'4.85"
$ws.Range("E44").Value = "  -2.33%  "

# Row 45
$ws.Range("D45").Value = "'0.577"
$ws.Range("E45").Value = "  -3.04%  "

# Row 46
$ws.Range("E46").Value = "  -1.01%  "

# Row 47
$ws.Range("D47").Value = "'119.76"
$ws.Range("E47").Value = "  -6.22%  "

# Row 48
$ws.Range("D48").Value = "'0.0481"
$ws.Range("E48").Value = "  -3.13%  "

# Row 49
$ws.Range("D49").Value = "'17.01"
$ws.Range("E49").Value = "  -4.82%  "

# Row 50
$ws.Range("D50").Value = "'0.0208"
$ws.Range("E50").Value = "  -3.75%  "

# Row 51
$ws.Range("D51").Value = "'16.43"
$ws.Range("E51").Value = "  -4.77%  "
